{"js": "// Rename the eviews_commands() call to eviews_rwalk() and change its\n// single string argument `\"wfcreate m 1990 +90\"` into a call with two\n// arguments: \"Eviewsr\", series = \"x y z\"  (matching the source-highlighted\n// R code block in the document).\n\n// 1) Rename the function-name token eviews_commands -> eviews_rwalk.\nconst fnResults = context.document.body.search(\"eviews_commands\", { matchCase: true });\nfnResults.load(\"items\");\nawait context.sync();\n\nif (fnResults.items.length === 0) {\n  throw new Error(\"Could not find 'eviews_commands' in the document.\");\n}\n\nfnResults.items[0].insertText(\"eviews_rwalk\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Replace the lone string-literal argument with the new first\n//    argument, then append the additional `series = \"x y z\"` argument,\n//    re-using the existing syntax-highlighting character styles\n//    (StringTok / NormalTok / AttributeTok) run-by-run.\nconst argResults = context.document.body.search('\"wfcreate m 1990 +90\"', { matchCase: true });\nargResults.load(\"items\");\nawait context.sync();\n\nif (argResults.items.length === 0) {\n  throw new Error('Could not find the \"wfcreate m 1990 +90\" argument in the document.');\n}\n\nconst firstArg = argResults.items[0];\nfirstArg.insertText('\"Eviewsr\"', Word.InsertLocation.replace);\nawait context.sync();\n\nlet cursor = firstArg.insertText(\",\", Word.InsertLocation.after);\ncursor.style = \"NormalTok\";\nawait context.sync();\n\ncursor = cursor.insertText(\"series =\", Word.InsertLocation.after);\ncursor.style = \"AttributeTok\";\nawait context.sync();\n\ncursor = cursor.insertText(\" \", Word.InsertLocation.after);\ncursor.style = \"NormalTok\";\nawait context.sync();\n\ncursor = cursor.insertText('\"x y z\"', Word.InsertLocation.after);\ncursor.style = \"StringTok\";\nawait context.sync();\n", "ps1": "# Rename the eviews_commands() call to eviews_rwalk() and change its\n# single string argument `\"wfcreate m 1990 +90\"` into a call with two\n# arguments: \"Eviewsr\", series = \"x y z\"  (matching the source-highlighted\n# R code block in the document).\n\n$d = $word.ActiveDocument\n\n# 1) Rename the function-name token eviews_commands -> eviews_rwalk.\n#    Find/Replace keeps the existing run formatting (FunctionTok style).\n$fnRange = $d.Content\n$fnFound = $fnRange.Find.Execute(\"eviews_commands\", $false, $false, $false, $false, $false, $true, 1, $false, \"eviews_rwalk\", 2)\nif (-not $fnFound) {\n    throw \"Could not find 'eviews_commands' in the document.\"\n}\n\n# 2) Replace the lone string-literal argument with the new first\n#    argument, then append the additional `series = \"x y z\"` argument,\n#    re-using the existing syntax-highlighting character styles\n#    (StringTok / NormalTok / AttributeTok) run-by-run.\n$argRange = $d.Content\n$argFound = $argRange.Find.Execute('\"wfcreate m 1990 +90\"', $true)\nif (-not $argFound) {\n    throw 'Could not find the \"wfcreate m 1990 +90\" argument in the document.'\n}\n\n$argRange.Text = '\"Eviewsr\"'\n$argRange.Collapse(0)\n\n$argRange.InsertAfter(\",\")\n$argRange.Style = \"NormalTok\"\n$argRange.Collapse(0)\n\n$argRange.InsertAfter(\"series =\")\n$argRange.Style = \"AttributeTok\"\n$argRange.Collapse(0)\n\n$argRange.InsertAfter(\" \")\n$argRange.Style = \"NormalTok\"\n$argRange.Collapse(0)\n\n$argRange.InsertAfter('\"x y z\"')\n$argRange.Style = \"StringTok\"\n"}
